$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'Al-Qaeda'
$ws.Range('B2').Value = 'Anti-communism,sfn,Sunni Islamism,Islamic fundamentalism,Factions:,Anti-Zionism,Collapsible list,loc,Anti-imperialism,Anti-Indian sentiment,87 * Sunni–Shia alliance sfn,Willsky-Ciollo,Salafi jihadism,Jihadism,Antisemitism,harvnb,Factions: * Pan-Islamism sfn,Anti-Americanism,Deobandism,Sunni–Shia alliance,Salafism,Anti-LGBT,Muslim unity,12 * Islamic fundamentalism * Anti-Americanism * Anti-communism * Anti-imperialism * Anti-Indian sentiment * Anti-LGBT * Antisemitism * Anti-Western imperialism * Anti-Zionism Plainlist,Anti-Western imperialism,Deobandi jihadism,Nbsp Plainlist,Gallagher,Wright,Sources:,Nbsp,Gunaratna,Bergen,Peter L.,Holy war,Inc.: Inside the Secret World of Osama bin Laden,New York: Free Press,2001.,pp. 70–71,Wahhabism,Qutbism,Ahl-i Hadith,officially denied),Introduction,pp. 12,87 * Qutbism * Jihadism * Muslim unity sfn,Introduction,pp. 12,87,Plainlist,Pan-Islamism'
$ws.Range('C2').Value = '-'
$ws.Range('D2').ClearContents() | Out-Null

$ws.Range('A3').Value = 'hy'
$ws.Range('B3').Value = 'Left-wing nationalism,Revolutionary socialism,Marxism-Leninism,Armenian nationalism'
$ws.Range('C3').Value = 'Left-wing to far-left'
$ws.Range('D3').ClearContents() | Out-Null

$ws.Range('A4').Value = 'Islamic Resistance Movement'
$ws.Range('B4').Value = 'sfn,Islamic fundamentalism,113 efn,Anti-Zionism,pp,156–57 sfn,66–67 efn,ubl,Palestinian nationalism,Militarism,Dalacoura,Gelvin,226 sfn,efn,465 sfn,Dunning,Litvak,66–67,66–67 sfn,156–57,Stepanova,Cheema,465 efn,Islamic nationalism,Islamism'
$ws.Range('C4').Value = '-'
$ws.Range('D4').Value = 'Yahya Sinwar,Khalil al-Hayya,Abu Omar Hassan,Yahya Sinwar,2024 targeted assassination of Muhammad Deif efn,Mohammed Deif Assassinated,'

$ws.Range('A5').Value = 'Lebanese Forces'
$ws.Range('B5').Value = 'Liberal conservatism,Lebanese nationalism,Christian democracy'
$ws.Range('C5').Value = 'Right-wing'
$ws.Range('D5').Value = 'Samir Geagea,Bachir Gemayel,Vice-president,'

$ws.Range('A6').Value = 'Lebanese Forces'
$ws.Range('B6').Value = 'Anti-communism,Federalism in Lebanon,Lebanese nationalism,Conservatism,Christian nationalism,Maronism,Anti Pan-Arabism'
$ws.Range('C6').Value = 'Right-wing to far-right'
$ws.Range('D6').ClearContents() | Out-Null

$ws.Range('A7').Value = 'Lebanese Front'
$ws.Range('B7').Value = 'Anti-Arabism,Anti-communism,Factions:,Phoenicianism,Anti-pan-Arabism,Anti-Palestinianism,Lebanese nationalism,Christian nationalism,Falangism'
$ws.Range('C7').Value = '-'
$ws.Range('D7').ClearContents() | Out-Null

$ws.Range('A8').Value = 'Palestine Liberation Organization'
$ws.Range('B8').Value = 'Pan-Arabism,Marxism,August 2024 * Factions:,Secularism,Factions:,ubl,Anti-Zionism,One-state solution,August 2024 * Baathism * Marxism,Palestinian nationalism,Anti-imperialism,Baathism,One-state solution * Anti-Zionism,August 2024,cn,Arab socialism,Arab nationalism'
$ws.Range('C8').Value = 'Left-wing'
$ws.Range('D8').Value = 'Mahmoud Abbas,'

$ws.Range('A9').Value = 'Palestine Liberation Organization'
$ws.Range('B9').Value = 'Pan-Arabism,Marxism,August 2024 * Factions:,Secularism,Factions:,ubl,Anti-Zionism,One-state solution,August 2024 * Baathism * Marxism,Palestinian nationalism,Anti-imperialism,Baathism,One-state solution * Anti-Zionism,August 2024,cn,Arab socialism,Arab nationalism'
$ws.Range('C9').Value = 'Left-wing'
$ws.Range('D9').Value = 'Mahmoud Abbas,'

$ws.Range('A10').Value = 'Lebanese Kataeb Party'
$ws.Range('B10').Value = 'Christian nationalism,Anti-communism,Social conservatism,Maronite politics,Lebanese nationalism,Christian democracy,Falangism'
$ws.Range('C10').Value = 'Right-wing,far-right,Centre-right'
$ws.Range('D10').Value = 'Samy Gemayel,Pierre Gemayel,'

$ws.Range('A11').Value = 'South Lebanon Army'
$ws.Range('B11').Value = 'Anti-communism,Muslim-Christian Unity,Secularism,Factions:,Zionism,Anti-Palestinianism,Multiconfessionalism,Maronite politics,Lebanese nationalism'
$ws.Range('C11').Value = '-'
$ws.Range('D11').Value = ','

$ws.Range('A12').Value = 'Syrian Social Nationalist Party'
$ws.Range('B12').Value = 'Antisemitism,Economic populism,Anti-communism,Social nationalism,ubl,collapsible list,Syrian irredentism,Syrian nationalism,Fascism'
$ws.Range('C12').Value = 'and right,In the past,the party or elements of its ideology or membership have been erred to as belonging to both the political left,sometimes being labeled far-right.,Syncretic efn,and right,sometimes being labeled far-right.'
$ws.Range('D12').Value = 'Rabie Banat,Antoun Saadeh,'

# Remove rows 13-19 content entirely (shrinks used range/dimension to D12)
$ws.Range('A13:D19').ClearContents() | Out-Null
